$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 values ---
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.3251496666666667
$ws.Range("N2").Value = 0.975449
$ws.Range("O2").Value = 0.07121046526627427
$ws.Range("P2").Value = 0.07121046526627427
$ws.Range("Q2").Value = 0.137297373097
$ws.Range("R2").Value = 1.235676357873
$ws.Range("S2").Value = 0.07121046526627427
$ws.Range("T2").Value = 0.07121046526627427

# --- Update existing row 3 values ---
$ws.Range("O3").Value = 0.2207208394324094
$ws.Range("P3").Value = 0.2207208394324094
$ws.Range("S3").Value = 0.2207208394324094
$ws.Range("T3").Value = 0.2207208394324094

# --- Update existing row 4 values ---
$ws.Range("M4").Value = 3.226895
$ws.Range("N4").Value = 9.680685
$ws.Range("O4").Value = 0.7067166842615477
$ws.Range("P4").Value = 0.7067166842615475
$ws.Range("Q4").Value = 1.362585455805
$ws.Range("R4").Value = 12.263269102245
$ws.Range("S4").Value = 0.7067166842615477
$ws.Range("T4").Value = 0.7067166842615475

# --- Add new row 5 ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt1"
$ws.Range("C5").Value = "Fzd3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.422259
$ws.Range("H5").Value = 1.266777
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.006173333333333333
$ws.Range("N5").Value = 0.01852
$ws.Range("O5").Value = 0.001352011039768762
$ws.Range("P5").Value = 0.001352011039768762
$ws.Range("Q5").Value = 0.00260674556
$ws.Range("R5").Value = 0.02346071004
$ws.Range("S5").Value = 0.001352011039768762
$ws.Range("T5").Value = 0.001352011039768762
